$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# Update the test-case label referenced by A3 (shared string reused in place)
$ws.Range("A3").Value = "MP.CPT.001.LEC"

# Clear the "last run" column D values while keeping formatting for D4:D8
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("D8").ClearContents()

# A6 had no real content; fully clear it (format + content) so the cell is removed
$ws.Range("A6").Clear()

# D9 was the counter start cell; fully clear it so the now-empty row disappears
$ws.Range("D9").Clear()

# Move the active selection to the counter start cell D8
$ws.Range("D8").Select()
